# Generated PowerShell-style Excel COM-interop edit script
# Applies the 'Updated symbol list' commit (Fri Dec 23 17:44:20 UTC 2022)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names, links, rank/volume labels): ---
# Assigning a non-numeric-looking string keeps Excel's COM layer from
# re-typing the cell as a number, so a direct .Value assignment is safe.
$textEdits = [ordered]@{
    "B10" = "One"
    "C10" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E10" = "9OneONEBestin24h"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "E11" = "10WazirXWRX"
    "B12" = "LiechtensteinCryptoassetsExchange"
    "C12" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "E12" = "11LiechtensteinCryptoassetsExchangeLCX"
    "B13" = "MandalaExchangeToken"
    "C13" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "E13" = "12MandalaExchangeTokenMDX"
    "B14" = "BitrueCoin"
    "C14" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "E14" = "13BitrueCoinBTR"
    "B15" = "MCDex"
    "C15" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "E15" = "14MCDexMCB"
    "B16" = "BitMartToken"
    "C16" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "E16" = "15BitMartTokenBMX"
    "B17" = "BitForexToken"
    "C17" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "E17" = "16BitForexTokenBF"
    "B18" = "CoinExToken"
    "C18" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "E18" = "17CoinExTokenCET"
    "E27" = "26UpBotsUBXT"
    "E48" = "47BOLOBOLOWorstin24h"
}
foreach ($addr in $textEdits.Keys) {
    $ws.Range($addr).Value = $textEdits[$addr]
}

# --- Numeric-looking text cells (prices): ---
# These must stay TEXT (to preserve exact formatting/trailing zeros, e.g.
# "0.01120" rather than the number 0.0112), so each is written with a
# leading apostrophe (Excel's literal-text prefix) and the cell style is
# then reset to Normal so no stray quote-prefix / text-format style sticks.
$priceEdits = [ordered]@{
    "D2" = "245.93"
    "D3" = "21.99"
    "D4" = "5.369"
    "D5" = "0.05852"
    "D8" = "0.8134"
    "D9" = "1.028"
    "D10" = "0.01120"
    "D11" = "0.1420"
    "D12" = "0.04319"
    "D13" = "0.07338"
    "D14" = "0.02987"
    "D15" = "4.165"
    "D16" = "0.09399"
    "D17" = "0.001587"
    "D18" = "0.04810"
    "D19" = "0.006063"
    "D20" = "0.004080"
    "D21" = "0.0009847"
    "D22" = "0.0001500"
    "D24" = "2.231"
    "D26" = "0.1268"
    "D41" = "0.006400"
    "D42" = "0.1073"
    "D43" = "0.003000"
    "D44" = "0.005071"
    "D48" = "0.07304"
}
foreach ($addr in $priceEdits.Keys) {
    $ws.Range($addr).Value = "'" + $priceEdits[$addr]
    $ws.Range($addr).Style = "Normal"
}
